$wb = $excel.ActiveWorkbook

$wsConsumption = $wb.Worksheets.Item("consumptionAssets")
$wsConversion  = $wb.Worksheets.Item("conversionAssets")

# --- consumptionAssets: new rows 6-9 ---
# Column B (asset name) filled top-to-bottom first
$wsConsumption.Cells.Item(6, 2).Value = "Industry_steel_electricity"
$wsConsumption.Cells.Item(7, 2).Value = "Industry_steel_heat"
$wsConsumption.Cells.Item(8, 2).Value = "Industry_other_electricity"
$wsConsumption.Cells.Item(9, 2).Value = "Industry_other_heat"

# Column C (alias) filled in this particular order
$wsConsumption.Cells.Item(6, 3).Value = "Industry_steel_electricity_demand"
$wsConsumption.Cells.Item(8, 3).Value = "Industry_other_electricity_demand"
$wsConsumption.Cells.Item(9, 3).Value = "Industry_other_heat_demand"
$wsConsumption.Cells.Item(7, 3).Value = "Industry_steel_heat_demand"

# Column E (energyAssetType)
$wsConsumption.Cells.Item(6, 5).Value = "ELECTRICITY_DEMAND"
$wsConsumption.Cells.Item(7, 5).Value = "HEAT_DEMAND"
$wsConsumption.Cells.Item(8, 5).Value = "ELECTRICITY_DEMAND"
$wsConsumption.Cells.Item(9, 5).Value = "HEAT_DEMAND"

# Column A (id), D (energyAssetCategory), F/G (values) - no new shared strings
$wsConsumption.Cells.Item(6, 1).Value = 5
$wsConsumption.Cells.Item(6, 4).Value = "CONSUMPTION"
$wsConsumption.Cells.Item(6, 6).Value = 1000000
$wsConsumption.Cells.Item(6, 7).Value = 0

$wsConsumption.Cells.Item(7, 1).Value = 6
$wsConsumption.Cells.Item(7, 4).Value = "CONSUMPTION"
$wsConsumption.Cells.Item(7, 6).Value = 0
$wsConsumption.Cells.Item(7, 7).Value = 1000000

$wsConsumption.Cells.Item(8, 1).Value = 7
$wsConsumption.Cells.Item(8, 4).Value = "CONSUMPTION"
$wsConsumption.Cells.Item(8, 6).Value = 1000000
$wsConsumption.Cells.Item(8, 7).Value = 0

$wsConsumption.Cells.Item(9, 1).Value = 8
$wsConsumption.Cells.Item(9, 4).Value = "CONSUMPTION"
$wsConsumption.Cells.Item(9, 6).Value = 0
$wsConsumption.Cells.Item(9, 7).Value = 1000000

$wsConsumption.Range("B10").Select()

# --- conversionAssets: new rows 10-11 ---
# Column B filled bottom row first, then top row
$wsConversion.Cells.Item(11, 2).Value = "Industrial_hydrogen_furnace"
$wsConversion.Cells.Item(10, 2).Value = "Industrial_methane_furnace"

# Column D (energyAssetType)
$wsConversion.Cells.Item(10, 4).Value = "METHANE_FURNACE"
$wsConversion.Cells.Item(11, 4).Value = "HYDROGEN_FURNACE"

# Remaining columns - no new shared strings
$wsConversion.Cells.Item(10, 1).Value = 9
$wsConversion.Cells.Item(10, 3).Value = "CONVERSION"
$wsConversion.Cells.Item(10, 5).Value = 0
$wsConversion.Cells.Item(10, 6).Value = 300
$wsConversion.Cells.Item(10, 7).Value = 0.99
$wsConversion.Cells.Item(10, 8).Value = 120

$wsConversion.Cells.Item(11, 1).Value = 10
$wsConversion.Cells.Item(11, 3).Value = "CONVERSION"
$wsConversion.Cells.Item(11, 5).Value = 0
$wsConversion.Cells.Item(11, 6).Value = 300
$wsConversion.Cells.Item(11, 7).Value = 0.99
$wsConversion.Cells.Item(11, 8).Value = 120

$wsConversion.Range("D12").Select()
